$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# Row 28 quantity changed from 2 to 1; the dependent "Total" formula in D28
# (shared formula B*C) recalculates automatically to 0.2, and the grand
# total in D54 (SUM(D4:D53)) recalculates accordingly.
$ws.Range("B28").Value = 1

# The active selection moved from B29 to B28.
$ws.Range("B28").Select()
